# Update EPEX Spot prices workbook with the latest daily data point.
$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": append a new date column (AN) with hourly prices ---
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Clone formatting from the neighbouring header cell (AM1) first so the new
# header lands in the same (reused) style slot as the rest of row 1, then set
# its text.
$wsSpot.Range("AM1").Copy() | Out-Null
$wsSpot.Range("AN1").PasteSpecial(-4122) | Out-Null
$wsSpot.Range("AN1").Value = "23-jul"

$spotValues = @(
    89.79000000000001,
    82.64,
    80.94,
    78.26000000000001,
    77.83,
    81.11,
    90.87,
    100.5,
    96.5,
    83.31,
    76.40000000000001,
    68.64,
    61.97,
    51.89,
    44.88,
    50.38,
    71.63,
    82.17,
    92.17,
    104.98,
    111.63,
    112,
    109.96,
    100.39
)

for ($i = 0; $i -lt $spotValues.Length; $i++) {
    $row = $i + 2
    $wsSpot.Cells.Item($row, 40).Value = $spotValues[$i]
}

# --- Sheet "Gaz": append the newest daily row ---
$wsGaz = $wb.Worksheets.Item("Gaz")
# Force the date to be stored as plain text (matching every other row in
# column A) instead of letting Excel auto-convert the "YYYY-MM-DD" string
# into a real date value.
$wsGaz.Range("A37").NumberFormat = "@"
$wsGaz.Range("A37").Value2 = "2025-07-21"
$wsGaz.Range("A37").Style = $wsGaz.Range("A36").Style
$wsGaz.Range("B37").Value = 32.6

# --- Sheet "CO2": append the newest daily row ---
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A37").NumberFormat = "@"
$wsCo2.Range("A37").Value2 = "2025-07-21"
$wsCo2.Range("A37").Style = $wsCo2.Range("A36").Style
$wsCo2.Range("B37").Value = 69.09999999999999
